$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("applicator")

$ws.Range("A65").Value = '''27/02/2018'
$ws.Range("B65").Value = '''80000571'
$ws.Range("C65").Value = '''3012'
$ws.Range("D65").Value = 'Гострини на розрізі контакту'
$ws.Range("E65").Value = '''19'

$ws.Range("A66").Value = '''04/03/2018'
$ws.Range("B66").Value = '''80000571'
$ws.Range("C66").Value = '''3012'
$ws.Range("D66").Value = 'Пошкодження поверхні контакту'
$ws.Range("E66").Value = '''20'

$ws.Range("A67").Value = '''04/03/2018'
$ws.Range("B67").Value = '''80000571'
$ws.Range("C67").Value = '''3012'
$ws.Range("D67").Value = 'Гострини на розрізі контакту'
$ws.Range("E67").Value = '''21'

$ws.Range("A68").Value = '''04/03/2018'
$ws.Range("B68").Value = '''80000571'
$ws.Range("C68").Value = '''3070'
$ws.Range("D68").Value = 'Гострини на розрізі контакту'
$ws.Range("E68").Value = '''22'

$ws.Range("A69").Value = '''04/03/2018'
$ws.Range("B69").Value = '''80000571'
$ws.Range("C69").Value = '''3070'
$ws.Range("D69").Value = 'Гострини на розрізі контакту'
$ws.Range("E69").Value = '''23'

$ws.Range("A70").Value = '''04/03/2018'
$ws.Range("B70").Value = '''80000571'
$ws.Range("C70").Value = '''3070'
$ws.Range("D70").Value = 'Гострини на розрізі контакту'
$ws.Range("E70").Value = '''24'

$ws.Range("A71").Value = '''04/03/2018'
$ws.Range("B71").Value = '''80000571'
$ws.Range("C71").Value = '''3070'
$ws.Range("D71").Value = 'Гострини на розрізі контакту'
$ws.Range("E71").Value = '''22'

$ws.Range("A72").Value = '''04/03/2018'
$ws.Range("B72").Value = '''80000571'
$ws.Range("C72").Value = '''3070'
$ws.Range("D72").Value = 'Гострини на розрізі контакту'
$ws.Range("E72").Value = '''23'

$ws.Range("A73").Value = '''04/03/2018'
$ws.Range("B73").Value = '''80000571'
$ws.Range("C73").Value = '''3070'
$ws.Range("D73").Value = 'Гострини на розрізі контакту'
$ws.Range("E73").Value = '''24'

$ws.Range("A74").Value = '''04/03/2018'
$ws.Range("B74").Value = '''80000571'
$ws.Range("C74").Value = '''3070'
$ws.Range("D74").Value = 'Гострини на розрізі контакту'
$ws.Range("E74").Value = '''25'

$ws.Range("A75").Value = '''04/03/2018'
$ws.Range("B75").Value = '''80000571'
$ws.Range("C75").Value = '''3070'
$ws.Range("D75").Value = 'Гострини на розрізі контакту'
$ws.Range("E75").Value = '''26'

$ws.Range("A76").Value = '''04/03/2018'
$ws.Range("B76").Value = '''80000571'
$ws.Range("C76").Value = '''3070'
$ws.Range("D76").Value = 'Гострини на розрізі контакту'
$ws.Range("E76").Value = '''27'

$ws.Range("A77").Value = '''04/03/2018'
$ws.Range("B77").Value = '''80000571'
$ws.Range("C77").Value = '''3070'
$ws.Range("D77").Value = 'Гострини на розрізі контакту'
$ws.Range("E77").Value = '''28'

$ws.Range("A78").Value = '''04/03/2018'
$ws.Range("B78").Value = '''80000571'
$ws.Range("C78").Value = '''3070'
$ws.Range("D78").Value = 'Гострини на розрізі контакту'
$ws.Range("E78").Value = '''29'

$ws.Range("A79").Value = '''04/03/2018'
$ws.Range("B79").Value = '''80000571'
$ws.Range("C79").Value = '''3070'
$ws.Range("D79").Value = 'Гострини на розрізі контакту'
$ws.Range("E79").Value = '''30'

$ws.Range("A80").Value = '''04/03/2018'
$ws.Range("B80").Value = '''80000571'
$ws.Range("C80").Value = '''3070'
$ws.Range("D80").Value = 'Гострини на розрізі контакту'
$ws.Range("E80").Value = '''31'

$ws.Range("A81").Value = '''04/03/2018'
$ws.Range("B81").Value = '''80000571'
$ws.Range("C81").Value = '''3070'
$ws.Range("D81").Value = 'Гострини на розрізі контакту'
$ws.Range("E81").Value = '''32'

$ws.Range("A82").Value = '''04/03/2018'
$ws.Range("B82").Value = '''80000571'
$ws.Range("C82").Value = '''3070'
$ws.Range("D82").Value = 'Гострини на розрізі контакту'
$ws.Range("E82").Value = '''33'

$ws.Range("A83").Value = '''04/03/2018'
$ws.Range("B83").Value = '''80000571'
$ws.Range("C83").Value = '''3070'
$ws.Range("D83").Value = 'Гострини на розрізі контакту'
$ws.Range("E83").Value = '''34'

$ws.Range("A84").Value = '''04/03/2018'
$ws.Range("B84").Value = '''80000571'
$ws.Range("C84").Value = '''3070'
$ws.Range("D84").Value = 'Гострини на розрізі контакту'
$ws.Range("E84").Value = '''35'

$ws.Range("A85").Value = '''04/03/2018'
$ws.Range("B85").Value = '''80000571'
$ws.Range("C85").Value = '''3070'
$ws.Range("D85").Value = 'Гострини на розрізі контакту'
$ws.Range("E85").Value = '''222'

$ws.Range("A86").Value = '''04/03/2018'
$ws.Range("B86").Value = '''80000571'
$ws.Range("C86").Value = '''3070'
$ws.Range("D86").Value = 'Гострини на розрізі контакту'
$ws.Range("E86").Value = '''223'

$ws.Range("A87").Value = '''04/03/2018'
$ws.Range("B87").Value = '''80000571'
$ws.Range("C87").Value = '''3070'
$ws.Range("D87").Value = 'Гострини на розрізі контакту'
$ws.Range("E87").Value = '''224'

$ws.Range("A88").Value = '''04/03/2018'
$ws.Range("B88").Value = '''80000571'
$ws.Range("C88").Value = '''3070'
$ws.Range("D88").Value = 'Гострини на розрізі контакту'
$ws.Range("E88").Value = '''225'

$ws.Range("A89").Value = '''04/03/2018'
$ws.Range("B89").Value = '''80000571'
$ws.Range("C89").Value = '''3070'
$ws.Range("D89").Value = 'Гострини на розрізі контакту'
$ws.Range("E89").Value = '''226'

$ws.Range("A90").Value = '''04/03/2018'
$ws.Range("B90").Value = '''80000571'
$ws.Range("C90").Value = '''3070'
$ws.Range("D90").Value = 'Гострини на розрізі контакту'
$ws.Range("E90").Value = '''227'

$ws.Range("A91").Value = '''04/03/2018'
$ws.Range("B91").Value = '''80000571'
$ws.Range("C91").Value = '''3070'
$ws.Range("D91").Value = 'Гострини на розрізі контакту'
$ws.Range("E91").Value = '''228'

$ws.Range("A92").Value = '''04/03/2018'
$ws.Range("B92").Value = '''80000571'
$ws.Range("C92").Value = '''3070'
$ws.Range("D92").Value = 'Гострини на розрізі контакту'
$ws.Range("E92").Value = '''229'

$ws.Range("A93").Value = '''04/03/2018'
$ws.Range("B93").Value = '''80000571'
$ws.Range("C93").Value = '''3070'
$ws.Range("D93").Value = 'Гострини на розрізі контакту'
$ws.Range("E93").Value = '''230'

$ws.Range("A94").Value = '''04/03/2018'
$ws.Range("B94").Value = '''80000571'
$ws.Range("C94").Value = '''3070'
$ws.Range("D94").Value = 'Гострини на розрізі контакту'
$ws.Range("E94").Value = '''233'

$ws.Range("A95").Value = '''04/03/2018'
$ws.Range("B95").Value = '''80000571'
$ws.Range("C95").Value = '''3070'
$ws.Range("D95").Value = 'Гострини на розрізі контакту'
$ws.Range("E95").Value = '''234'

$ws.Range("A96").Value = '''04/03/2018'
$ws.Range("B96").Value = '''80000571'
$ws.Range("C96").Value = '''3070'
$ws.Range("D96").Value = 'Гострини на розрізі контакту'
$ws.Range("E96").Value = '''235'

$ws.Range("A97").Value = '''04/03/2018'
$ws.Range("B97").Value = '''80000571'
$ws.Range("C97").Value = '''3070'
$ws.Range("D97").Value = 'Гострини на розрізі контакту'
$ws.Range("E97").Value = '''236'

$ws.Range("A98").Value = '''04/03/2018'
$ws.Range("B98").Value = '''80000571'
$ws.Range("C98").Value = '''3070'
$ws.Range("D98").Value = 'Гострини на розрізі контакту'
$ws.Range("E98").Value = '''237'

$ws.Range("A99").Value = '''04/03/2018'
$ws.Range("B99").Value = '''80000571'
$ws.Range("C99").Value = '''3070'
$ws.Range("D99").Value = 'Не симетричне / не відповідне закриття ядра'
$ws.Range("E99").Value = '''238'

$ws.Range("A100").Value = '''04/03/2018'
$ws.Range("B100").Value = '''80000571'
$ws.Range("C100").Value = '''3070'
$ws.Range("D100").Value = 'Не симетричне / не відповідне закриття ядра'
$ws.Range("E100").Value = '''238'

$ws.Range("A101").Value = '''05/03/2018'
$ws.Range("B101").Value = '''80000571'
$ws.Range("C101").Value = '''3070'
$ws.Range("D101").Value = 'Гострини на розрізі контакту'
$ws.Range("E101").Value = '''240'

$ws.Range("A102").Value = '''05/03/2018'
$ws.Range("B102").Value = '''80000571'
$ws.Range("C102").Value = '''3070'
$ws.Range("D102").Value = 'Гострини на розрізі контакту'
$ws.Range("E102").Value = '''241'

$ws.Range("A103").Value = '''05/03/2018'
$ws.Range("B103").Value = '''80000571'
$ws.Range("C103").Value = '''3070'
$ws.Range("D103").Value = 'Гострини на розрізі контакту'
$ws.Range("E103").Value = '''241'

$ws.Range("A104").Value = '''05/03/2018'
$ws.Range("B104").Value = '''80000571'
$ws.Range("C104").Value = '''3070'
$ws.Range("D104").Value = 'Гострини на розрізі контакту'
$ws.Range("E104").Value = '''241'

$ws.Range("A105").Value = '''05/03/2018'
$ws.Range("B105").Value = '''80000571'
$ws.Range("C105").Value = '''3070'
$ws.Range("D105").Value = 'Гострини на розрізі контакту'
$ws.Range("E105").Value = '''242'

$ws.Range("A106").Value = '''05/03/2018'
$ws.Range("B106").Value = '''80000571'
$ws.Range("C106").Value = '''3070'
$ws.Range("D106").Value = 'Пошкодження поверхні контакту'
$ws.Range("E106").Value = '''243'

$ws.Range("A107").Value = '''05/03/2018'
$ws.Range("B107").Value = '''80000571'
$ws.Range("C107").Value = '''3070'
$ws.Range("D107").Value = 'Гострини на розрізі контакту'
$ws.Range("E107").Value = '''244'

$ws.Range("A108").Value = '**'
